$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.572.40"
$ws.Range("E2").Value = "  +3.94%  "
$ws.Range("D3").Value = "3.141.53"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "395.83"
$ws.Range("E5").Value = "  +2.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.65"
$ws.Range("E6").Value = "  +6.40%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("E9").Value = "  +4.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.80"
$ws.Range("E10").Value = "  +5.58%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").Value = "3.644.03"
$ws.Range("E13").Value = "  +2.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.00"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.99"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("E16").Value = "  +8.41%  "
$ws.Range("D17").Value = "3.148.44"
$ws.Range("E17").Value = "  +3.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.48"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("D19").Value = "53.420.14"
$ws.Range("E19").Value = "  +3.53%  "
$ws.Range("E20").Value = "  +3.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.74"
$ws.Range("E21").Value = "  +2.36%  "
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.81"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.43"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.95"
$ws.Range("E26").Value = "  -3.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.47"
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.34"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.92"
$ws.Range("E32").Value = "  +6.37%  "
$ws.Range("E33").Value = "  +11.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "36.91"
$ws.Range("E34").Value = "  +6.20%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.36"
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.64"
$ws.Range("E37").Value = "  +9.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.78"
$ws.Range("E39").Value = "  +8.67%  "
$ws.Range("E40").Value = "  +9.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.288"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.07"
$ws.Range("E44").Value = "  +3.78%  "
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.12"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.08"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "2.073.11"
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("E50").Value = "  +5.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0499"
$ws.Range("E51").Value = "  +14.65%  "
